$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-90
# from serial date 45175 (2023-09-06) to 45177 (2023-09-08)
$ws.Range("C2:C90").Value = 45177
